# Auto-generated Excel COM-interop script to apply Durandal_Profits.xlsx market-data update
$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1496.4375
$ws.Range("I62").Value = 1637.1875
$ws.Range("J62").Value = 1355.6875
$ws.Range("K62").Value = 1637.1875
$ws.Range("L62").Value = 1355.6875
$ws.Range("M62").Value = -1013.1875
$ws.Range("N62").Value = -2603.6875

# Hunk 1: sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1496.4375
$ws.Range("I65").Value = 1637.1875
$ws.Range("J65").Value = 1355.6875
$ws.Range("K65").Value = 8185.9375
$ws.Range("L65").Value = 6778.4375
$ws.Range("M65").Value = -5065.9375
$ws.Range("N65").Value = -13018.4375

# Hunk 2: sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 637.37933
$ws.Range("I107").Value = 769.0454999999999
$ws.Range("J107").Value = 223.57143
$ws.Range("K107").Value = 769.0454999999999
$ws.Range("L107").Value = 223.57143
$ws.Range("M107").Value = 1150.9545
$ws.Range("N107").Value = -4063.57143

# Hunk 3: sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7976.7646
$ws.Range("I116").Value = 15900.714
$ws.Range("J116").Value = 2430
$ws.Range("K116").Value = 15900.714
$ws.Range("L116").Value = 2430
$ws.Range("M116").Value = -12458.714
$ws.Range("N116").Value = -9314

# Hunk 4: sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1083.5238
$ws.Range("I137").Value = 871.1
$ws.Range("J137").Value = 1276.6364
$ws.Range("K137").Value = 2613.3
$ws.Range("L137").Value = 3829.9092
$ws.Range("M137").Value = -63.30000000000018
$ws.Range("N137").Value = -8929.9092

# Hunk 5: sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2771.6172
$ws.Range("I138").Value = 1210.8334
$ws.Range("J138").Value = 3689.7256
$ws.Range("K138").Value = 3632.5002
$ws.Range("L138").Value = 11069.1768
$ws.Range("M138").Value = 1507.4998
$ws.Range("N138").Value = -21349.1768

# Hunk 6: sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 415220.25
$ws.Range("I32").Value = 3975.1526
$ws.Range("J32").Value = 2620989.5
$ws.Range("K32").Value = 3975.1526
$ws.Range("L32").Value = 2620989.5
$ws.Range("M32").Value = -3688.1526
$ws.Range("N32").Value = -2621563.5

# Hunk 7: sheet BSM, row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 10306.7
$ws.Range("I35").Value = 2067
$ws.Range("J35").Value = 11222.223
$ws.Range("K35").Value = 2067
$ws.Range("L35").Value = 11222.223
$ws.Range("M35").Value = -1757
$ws.Range("N35").Value = -11842.223

# Hunk 8: sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2727.25
$ws.Range("I134").Value = 2721.7144
$ws.Range("J134").Value = 2766
$ws.Range("K134").Value = 8165.1432
$ws.Range("L134").Value = 8298
$ws.Range("M134").Value = -5630.1432
$ws.Range("N134").Value = -13368

# Hunk 9: sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6383.615
$ws.Range("I31").Value = 6908.2383
$ws.Range("J31").Value = 4180.2
$ws.Range("K31").Value = 6908.2383
$ws.Range("L31").Value = 4180.2
$ws.Range("M31").Value = -6613.2383
$ws.Range("N31").Value = -4770.2

# Hunk 10: sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6383.615
$ws.Range("I34").Value = 6908.2383
$ws.Range("J34").Value = 4180.2
$ws.Range("K34").Value = 6908.2383
$ws.Range("L34").Value = 4180.2
$ws.Range("M34").Value = -6706.2383
$ws.Range("N34").Value = -4584.2

# Hunk 11: sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3081.3928
$ws.Range("I58").Value = 1149.5
$ws.Range("J58").Value = 3854.15
$ws.Range("K58").Value = 1149.5
$ws.Range("L58").Value = 3854.15
$ws.Range("M58").Value = -946.5
$ws.Range("N58").Value = -4260.15

# Hunk 12: sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1213.3636
$ws.Range("I122").Value = 730
$ws.Range("J122").Value = 1394.625
$ws.Range("K122").Value = 2190
$ws.Range("L122").Value = 4183.875
$ws.Range("M122").Value = 260
$ws.Range("N122").Value = -9083.875

# Hunk 13: sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 45900.477
$ws.Range("I132").Value = 1668.625
$ws.Range("K132").Value = 5005.875
$ws.Range("M132").Value = -2475.875

# Hunk 14: sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3081.3928
$ws.Range("I136").Value = 1149.5
$ws.Range("J136").Value = 3854.15
$ws.Range("K136").Value = 3448.5
$ws.Range("L136").Value = 11562.45
$ws.Range("M136").Value = -898.5
$ws.Range("N136").Value = -16662.45

# Hunk 15: sheet CUL, row 88
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 4705.3335
$ws.Range("J88").Value = 4705.3335
$ws.Range("L88").Value = 14116.0005
$ws.Range("N88").Value = -14972.0005

# Hunk 16: sheet CUL, row 91
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 4705.3335
$ws.Range("J91").Value = 4705.3335
$ws.Range("L91").Value = 14116.0005
$ws.Range("N91").Value = -17080.0005

# Hunk 17: sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2237.5
$ws.Range("I102").Value = 1716.8334
$ws.Range("J102").Value = 3799.5
$ws.Range("K102").Value = 1716.8334
$ws.Range("L102").Value = 3799.5
$ws.Range("M102").Value = -94.83339999999998
$ws.Range("N102").Value = -7043.5

# Hunk 18: sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1653.1538
$ws.Range("I113").Value = 1459.1
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 1459.1
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 710.9000000000001
$ws.Range("N113").Value = -6640

# Hunk 19: sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 35092.9
$ws.Range("I132").Value = 1283.8572
$ws.Range("J132").Value = 113980.664
$ws.Range("K132").Value = 3851.5716
$ws.Range("L132").Value = 341941.992
$ws.Range("M132").Value = -1321.5716
$ws.Range("N132").Value = -347001.992

# Hunk 20: sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1333.5714
$ws.Range("I7").Value = 1405.4166
$ws.Range("J7").Value = 902.5
$ws.Range("K7").Value = 1405.4166
$ws.Range("L7").Value = 902.5
$ws.Range("M7").Value = -1293.4166
$ws.Range("N7").Value = -1126.5

# Hunk 21: sheet LTW, row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 4000
$ws.Range("J34").Value = 4000
$ws.Range("L34").Value = 4000
$ws.Range("N34").Value = -4344

# Hunk 22: sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1333.5714
$ws.Range("I126").Value = 1405.4166
$ws.Range("J126").Value = 902.5
$ws.Range("K126").Value = 4216.2498
$ws.Range("L126").Value = 2707.5
$ws.Range("M126").Value = -1746.2498
$ws.Range("N126").Value = -7647.5

# Hunk 23: sheet LTW, row 131
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 10132.286
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 10132.286
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 10132.286
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -20212.286

# Hunk 24: sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1388.6
$ws.Range("I122").Value = 1384.2106
$ws.Range("J122").Value = 1402.5
$ws.Range("K122").Value = 4152.6318
$ws.Range("L122").Value = 4207.5
$ws.Range("M122").Value = -1702.6318
$ws.Range("N122").Value = -9107.5

# Hunk 25: sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19407.611
$ws.Range("I136").Value = 31093.969
$ws.Range("J136").Value = 1043.3334
$ws.Range("K136").Value = 93281.90700000001
$ws.Range("L136").Value = 3130.0002
$ws.Range("N136").Value = -8230.0002
